$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "45.289.07"
Set-TextValue $ws.Range("E2") "  +3.87%  "
Set-TextValue $ws.Range("D3") "2.422.71"
Set-TextValue $ws.Range("E3") "  +0.03%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.06%  "
Set-TextValue $ws.Range("D5") "318.19"
Set-TextValue $ws.Range("E5") "  +3.81%  "
Set-TextValue $ws.Range("D6") "102.37"
Set-TextValue $ws.Range("E6") "  +5.52%  "
Set-TextValue $ws.Range("D7") "0.515"
Set-TextValue $ws.Range("E7") "  +1.19%  "
Set-TextValue $ws.Range("E8") "  -0.06%  "
Set-TextValue $ws.Range("D9") "0.529"
Set-TextValue $ws.Range("E9") "  +7.73%  "
Set-TextValue $ws.Range("D10") "35.54"
Set-TextValue $ws.Range("E10") "  +1.35%  "
Set-TextValue $ws.Range("D11") "0.0803"
Set-TextValue $ws.Range("E11") "  +0.74%  "
Set-TextValue $ws.Range("D13") "18.17"
Set-TextValue $ws.Range("E13") "  -1.58%  "
Set-TextValue $ws.Range("D14") "7.02"
Set-TextValue $ws.Range("E14") "  +1.85%  "
Set-TextValue $ws.Range("D15") "2.804.11"
Set-TextValue $ws.Range("E15") "  +0.30%  "
Set-TextValue $ws.Range("D16") "2.426.13"
Set-TextValue $ws.Range("E16") "  +0.88%  "
Set-TextValue $ws.Range("E17") "  +1.75%  "
Set-TextValue $ws.Range("D18") "45.195.39"
Set-TextValue $ws.Range("E18") "  +3.60%  "
Set-TextValue $ws.Range("D19") "12.21"
Set-TextValue $ws.Range("E19") "  +1.28%  "
Set-TextValue $ws.Range("E20") "  -1.24%  "
Set-TextValue $ws.Range("E21") "  +2.16%  "
Set-TextValue $ws.Range("D22") "68.76"
Set-TextValue $ws.Range("E22") "  +0.96%  "
Set-TextValue $ws.Range("D23") "243.77"
Set-TextValue $ws.Range("E23") "  +2.55%  "
Set-TextValue $ws.Range("E24") "  +0.17%  "
Set-TextValue $ws.Range("E25") "  +1.99%  "
Set-TextValue $ws.Range("E26") "  -0.04%  "
Set-TextValue $ws.Range("D27") "25.50"
Set-TextValue $ws.Range("E27") "  +2.03%  "
Set-TextValue $ws.Range("D28") "9.57"
Set-TextValue $ws.Range("E28") "  +1.34%  "
Set-TextValue $ws.Range("E29") "  -11.90%  "
Set-TextValue $ws.Range("E30") "  +2.00%  "
Set-TextValue $ws.Range("D31") "32.79"
Set-TextValue $ws.Range("E31") "  +1.52%  "
Set-TextValue $ws.Range("D32") "0.126"
Set-TextValue $ws.Range("E32") "  +5.59%  "
Set-TextValue $ws.Range("D33") "20.17"
Set-TextValue $ws.Range("E33") "  +9.45%  "
Set-TextValue $ws.Range("D34") "5.20"
Set-TextValue $ws.Range("E34") "  +1.28%  "
Set-TextValue $ws.Range("E35") "  +0.23%  "
Set-TextValue $ws.Range("E36") "  +1.99%  "
Set-TextValue $ws.Range("E37") "  -1.46%  "
Set-TextValue $ws.Range("D38") "4.44"
Set-TextValue $ws.Range("D39") "126.60"
Set-TextValue $ws.Range("E39") "  -2.87%  "
Set-TextValue $ws.Range("D40") "2.85"
Set-TextValue $ws.Range("E40") "  -2.57%  "
Set-TextValue $ws.Range("E41") "  -2.40%  "
Set-TextValue $ws.Range("D42") "0.109"
Set-TextValue $ws.Range("E42") "  +0.86%  "
Set-TextValue $ws.Range("D43") "20.42"
Set-TextValue $ws.Range("E43") "  -2.39%  "
Set-TextValue $ws.Range("D44") "0.0289"
Set-TextValue $ws.Range("E44") "  +2.38%  "
Set-TextValue $ws.Range("D45") "1.931.14"
Set-TextValue $ws.Range("E45") "  -0.70%  "
Set-TextValue $ws.Range("E46") "  -2.93%  "
Set-TextValue $ws.Range("D47") "2.92"
Set-TextValue $ws.Range("E47") "  +3.12%  "
Set-TextValue $ws.Range("E48") "  +15.86%  "
Set-TextValue $ws.Range("D49") "9.10"
Set-TextValue $ws.Range("E49") "  -2.20%  "
Set-TextValue $ws.Range("D50") "76.50"
Set-TextValue $ws.Range("E50") "  +5.69%  "
Set-TextValue $ws.Range("D51") "53.85"
Set-TextValue $ws.Range("E51") "  +2.35%  "
